# "Fixed import of cashflows"
# The Category*/Sub Category*/Investment Domicile* columns (H, I, J) were
# dropped from the cashflow import template. The former Tag/Instrument
# columns (K, L) shift left into H/I, the "Instrument" header becomes
# "Instrument *", and the help-comments on the (new) H1/I1 cells are
# rewritten to reflect that Tag is optional and Instrument is mandatory.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old J1 comment ("Investment Domicile *") before the column
# holding it disappears.
$null = $ws.Range("J1").Comment.Delete()

# Drop columns H:J (Category *, Sub Category *, Investment Domicile *).
# This shifts old columns K (Tag) and L (Instrument) left to H and I.
$ws.Range("H1:J5").EntireColumn.Delete()

# Header text tweak: "Instrument" -> "Instrument *"
$ws.Range("I1").Value = "Instrument *"

# Update the help comments now sitting on H1 (Tag) and I1 (Instrument *).
$null = $ws.Range("H1").Comment.Text("Author:" + [char]10 + "-Optional")
$null = $ws.Range("I1").Comment.Text("Author:" + [char]10 + "-Mandatory")

# Match the saved selection/active cell.
$null = $ws.Range("I3").Select()

# Re-apply explicit column widths for the columns whose content changed
# (D kept its own width tweak upstream; F/G shrank back to fit; H/I are
# now the narrower Tag/Instrument columns).
$ws.Columns.Item(4).ColumnWidth = 10.084
$ws.Columns.Item(6).ColumnWidth = 9.084
$ws.Columns.Item(7).ColumnWidth = 8.75
$ws.Columns.Item(8).ColumnWidth = 6.084
$ws.Columns.Item(9).ColumnWidth = 13.084
